{"js": "// Replace the three-digit-by-one-digit multiplication answers in the\n// worksheet table with the new set of problems/answers, matching the\n// target revision exactly (old text -> new text), one-to-one.\nconst replacements = [\n  [\"117\u00d78=936\", \"647\u00d76=3882\"],\n  [\"853\u00d78=6824\", \"221\u00d76=1326\"],\n  [\"527\u00d73=1581\", \"562\u00d73=1686\"],\n  [\"972\u00d76=5832\", \"400\u00d72=800\"],\n  [\"135\u00d79=1215\", \"260\u00d76=1560\"],\n  [\"347\u00d73=1041\", \"514\u00d73=1542\"],\n  [\"465\u00d78=3720\", \"927\u00d72=1854\"],\n  [\"510\u00d77=3570\", \"296\u00d74=1184\"],\n  [\"405\u00d78=3240\", \"728\u00d77=5096\"],\n  [\"916\u00d79=8244\", \"978\u00d72=1956\"],\n  [\"107\u00d72=214\", \"173\u00d79=1557\"],\n  [\"977\u00d79=8793\", \"880\u00d75=4400\"],\n  [\"830\u00d76=4980\", \"435\u00d77=3045\"],\n  [\"584\u00d76=3504\", \"294\u00d72=588\"],\n  [\"951\u00d76=5706\", \"493\u00d72=986\"],\n  [\"354\u00d79=3186\", \"624\u00d73=1872\"],\n  [\"935\u00d78=7480\", \"530\u00d77=3710\"],\n  [\"401\u00d77=2807\", \"478\u00d74=1912\"],\n  [\"667\u00d73=2001\", \"638\u00d77=4466\"],\n  [\"654\u00d72=1308\", \"617\u00d74=2468\"],\n  [\"473\u00d75=2365\", \"649\u00d75=3245\"],\n  [\"373\u00d75=1865\", \"961\u00d76=5766\"],\n  [\"215\u00d72=430\", \"374\u00d73=1122\"],\n  [\"487\u00d76=2922\", \"135\u00d75=675\"],\n  [\"359\u00d78=2872\", \"863\u00d76=5178\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Pairs of (old text, new text) mirroring the target revision, applied\n# as a sequence of exact whole-match Find/Replace operations so each\n# answer cell is updated independently.\n$pairs = @(\n    @(\"117\u00d78=936\", \"647\u00d76=3882\"),\n    @(\"853\u00d78=6824\", \"221\u00d76=1326\"),\n    @(\"527\u00d73=1581\", \"562\u00d73=1686\"),\n    @(\"972\u00d76=5832\", \"400\u00d72=800\"),\n    @(\"135\u00d79=1215\", \"260\u00d76=1560\"),\n    @(\"347\u00d73=1041\", \"514\u00d73=1542\"),\n    @(\"465\u00d78=3720\", \"927\u00d72=1854\"),\n    @(\"510\u00d77=3570\", \"296\u00d74=1184\"),\n    @(\"405\u00d78=3240\", \"728\u00d77=5096\"),\n    @(\"916\u00d79=8244\", \"978\u00d72=1956\"),\n    @(\"107\u00d72=214\", \"173\u00d79=1557\"),\n    @(\"977\u00d79=8793\", \"880\u00d75=4400\"),\n    @(\"830\u00d76=4980\", \"435\u00d77=3045\"),\n    @(\"584\u00d76=3504\", \"294\u00d72=588\"),\n    @(\"951\u00d76=5706\", \"493\u00d72=986\"),\n    @(\"354\u00d79=3186\", \"624\u00d73=1872\"),\n    @(\"935\u00d78=7480\", \"530\u00d77=3710\"),\n    @(\"401\u00d77=2807\", \"478\u00d74=1912\"),\n    @(\"667\u00d73=2001\", \"638\u00d77=4466\"),\n    @(\"654\u00d72=1308\", \"617\u00d74=2468\"),\n    @(\"473\u00d75=2365\", \"649\u00d75=3245\"),\n    @(\"373\u00d75=1865\", \"961\u00d76=5766\"),\n    @(\"215\u00d72=430\", \"374\u00d73=1122\"),\n    @(\"487\u00d76=2922\", \"135\u00d75=675\"),\n    @(\"359\u00d78=2872\", \"863\u00d76=5178\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
